# Update cryptos list values per upstream scrape (Mon Jun 26 23:34:04 UTC 2023)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.247.49"
$ws.Range("E2").Value = "  -0.73%  "

$ws.Range("D3").Value = "1.857.24"

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Formula = "'236.38"
$ws.Range("E5").Value = "  -1.06%  "

$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").Formula = "'0.4788"
$ws.Range("E7").Value = "  -2.37%  "

$ws.Range("D8").Formula = "'0.2801"

$ws.Range("E9").Value = "  -3.19%  "

$ws.Range("D10").Value = "1.849.76"
$ws.Range("E10").Value = "  -2.54%  "

$ws.Range("D11").Formula = "'0.07390"
$ws.Range("E11").Value = "  +0.82%  "

$ws.Range("D12").Formula = "'16.22"
$ws.Range("E12").Value = "  -4.26%  "

$ws.Range("D13").Formula = "'5.015"
$ws.Range("E13").Value = "  -3.15%  "

$ws.Range("D14").Formula = "'87.09"
$ws.Range("E14").Value = "  -1.10%  "

$ws.Range("D15").Formula = "'0.6443"
$ws.Range("E15").Value = "  -3.48%  "

$ws.Range("D16").Value = "30.184.98"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").Formula = "'1.001"
$ws.Range("E17").Value = "  +0.03%  "

$ws.Range("D18").Formula = "'13.12"
$ws.Range("E18").Value = "  -2.19%  "

$ws.Range("D19").Formula = "'0.000007553"
$ws.Range("E19").Value = "  -4.05%  "

$ws.Range("D20").Value = "2.101.08"
$ws.Range("E20").Value = "  -2.05%  "

$ws.Range("D21").Formula = "'222.26"
$ws.Range("E21").Value = "  +14.07%  "

$ws.Range("D22").Formula = "'1.001"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Formula = "'5.272"
$ws.Range("E23").Value = "  -2.88%  "

$ws.Range("D24").Formula = "'6.067"
$ws.Range("E24").Value = "  -1.05%  "

$ws.Range("D25").Formula = "'9.189"
$ws.Range("E25").Value = "  -3.22%  "

$ws.Range("D26").Formula = "'163.64"
$ws.Range("E26").Value = "  +0.63%  "

$ws.Range("D27").Formula = "'18.54"
$ws.Range("E27").Value = "  +0.92%  "

$ws.Range("D28").Formula = "'1.927"
$ws.Range("E28").Value = "  -0.60%  "

$ws.Range("E29").Value = "  -3.59%  "

$ws.Range("D30").Formula = "'0.09190"
$ws.Range("E30").Value = "  +0.49%  "

$ws.Range("D31").Formula = "'4.234"
$ws.Range("E31").Value = "  -2.22%  "

$ws.Range("D32").Formula = "'3.961"
$ws.Range("E32").Value = "  -3.80%  "

$ws.Range("D33").Formula = "'0.04953"
$ws.Range("E33").Value = "  -3.97%  "

$ws.Range("D34").Formula = "'1.140"
$ws.Range("E34").Value = "  +3.11%  "

$ws.Range("D35").Formula = "'0.7226"
$ws.Range("E35").Value = "  -2.13%  "

$ws.Range("D36").Formula = "'2.687"
$ws.Range("E36").Value = "  -1.52%  "

$ws.Range("E37").Value = "  -1.07%  "

$ws.Range("D38").Formula = "'2.599"
$ws.Range("E38").Value = "  -2.75%  "

$ws.Range("D39").Formula = "'0.8976"
$ws.Range("E39").Value = "  -2.95%  "

$ws.Range("D40").Formula = "'2.036"
$ws.Range("E40").Value = "  -1.61%  "

$ws.Range("D41").Formula = "'5.907"
$ws.Range("E41").Value = "  +0.16%  "

$ws.Range("D42").Formula = "'105.94"
$ws.Range("E42").Value = "  -0.81%  "

$ws.Range("E43").Value = "  +0.51%  "

$ws.Range("D44").Formula = "'0.4237"
$ws.Range("E44").Value = "  -3.83%  "

$ws.Range("D45").Formula = "'7.263"
$ws.Range("E45").Value = "  -3.98%  "

$ws.Range("D46").Formula = "'0.1296"
$ws.Range("E46").Value = "  -5.40%  "

$ws.Range("D47").Formula = "'63.44"
$ws.Range("E47").Value = "  -8.17%  "

$ws.Range("D48").Formula = "'1.489"
$ws.Range("E48").Value = "  +6.07%  "

$ws.Range("D49").Formula = "'8.713"
$ws.Range("E49").Value = "  -3.31%  "

$ws.Range("D50").Formula = "'33.73"
$ws.Range("E50").Value = "  -3.55%  "

$ws.Range("E51").Value = "  -3.37%  "
